$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
}

Set-TextValue "D2" "66.380.50"
Set-TextValue "E2" "  +2.66%  "
Set-TextValue "D3" "3.243.57"
Set-TextValue "E3" "  +5.10%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "577.55"
Set-TextValue "E5" "  +2.16%  "
Set-TextValue "D6" "155.04"
Set-TextValue "E6" "  +8.02%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "3.235.09"
Set-TextValue "E8" "  +5.16%  "
Set-TextValue "E9" "  +4.16%  "
Set-TextValue "D10" "7.06"
Set-TextValue "E10" "  +10.18%  "
Set-TextValue "E11" "  +4.42%  "
Set-TextValue "D12" "0.490"
Set-TextValue "E12" "  +3.97%  "
Set-TextValue "D13" "37.91"
Set-TextValue "E13" "  +4.94%  "
Set-TextValue "E14" "  +4.07%  "
Set-TextValue "D15" "3.761.59"
Set-TextValue "E15" "  +5.06%  "
Set-TextValue "D16" "565.12"
Set-TextValue "E16" "  +13.48%  "
Set-TextValue "D17" "66.463.26"
Set-TextValue "E17" "  +2.74%  "
Set-TextValue "D18" "3.244.62"
Set-TextValue "E18" "  +5.05%  "
Set-TextValue "E19" "  +3.03%  "
Set-TextValue "D20" "7.13"
Set-TextValue "E20" "  +5.72%  "
Set-TextValue "D21" "14.45"
Set-TextValue "E21" "  +4.21%  "
Set-TextValue "D22" "0.747"
Set-TextValue "E22" "  +7.42%  "
Set-TextValue "D23" "7.90"
Set-TextValue "D24" "13.63"
Set-TextValue "E24" "  +6.30%  "
Set-TextValue "D25" "82.34"
Set-TextValue "E25" "  +3.80%  "
Set-TextValue "E26" "  -0.21%  "
Set-TextValue "D27" "9.44"
Set-TextValue "E27" "  +17.85%  "
Set-TextValue "E28" "  +5.75%  "
Set-TextValue "E29" "  +7.50%  "
Set-TextValue "D30" "28.12"
Set-TextValue "E30" "  +6.00%  "
Set-TextValue "E31" "  +2.05%  "
Set-TextValue "E32" "  +0.04%  "
Set-TextValue "E33" "  +4.00%  "
Set-TextValue "D34" "567.07"
Set-TextValue "E34" "  +8.32%  "
Set-TextValue "D35" "5.78"
Set-TextValue "E35" "  +3.87%  "
Set-TextValue "D36" "6.43"
Set-TextValue "E36" "  +6.60%  "
Set-TextValue "B37" "OKB"
Set-TextValue "C37" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D37" "55.96"
Set-TextValue "E37" "  +4.50%  "
Set-TextValue "B38" "VeChain"
Set-TextValue "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D38" "0.0462"
Set-TextValue "E38" "  +12.90%  "
Set-TextValue "D39" "0.0869"
Set-TextValue "E39" "  +7.43%  "
Set-TextValue "D40" "3.05"
Set-TextValue "E40" "  +13.45%  "
Set-TextValue "D41" "0.127"
Set-TextValue "E41" "  +5.43%  "
Set-TextValue "D42" "3.147.16"
Set-TextValue "E42" "  +6.87%  "
Set-TextValue "D43" "8.66"
Set-TextValue "E43" "  +2.34%  "
Set-TextValue "E44" "  +10.19%  "
Set-TextValue "D45" "2.34"
Set-TextValue "E45" "  +6.76%  "
Set-TextValue "D46" "26.80"
Set-TextValue "E46" "  +4.68%  "
Set-TextValue "D47" "0.0₃0564"
Set-TextValue "E47" "  +3.24%  "
Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  +0.10%  "
Set-TextValue "E49" "  +3.57%  "
Set-TextValue "E50" "  +8.11%  "
Set-TextValue "D51" "122.61"
Set-TextValue "E51" "  +1.78%  "
